$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Docentes responsáveis:" value row (row 13) entirely -
# this shifts every row below it up by one.
$ws.Rows(13).Delete()

# Update the remaining cells with their new content.
$ws.Range("B10:C10").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15:C15").Value = "01/01/2018"
$ws.Range("B18:C18").Value = "101761 - Arnaldo Márcio Ramalho Prata"
